# Module - Investment_Opening Status - Completed
# Adds a new "Debenture_series_setting" module-status sheet, fills in the
# InvestmentOpn sheet's completed test-run data, and updates the active
# tab / selections accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Add the new "Debenture_series_setting" sheet after Centre_Registration
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Debenture_series_setting"
$centre = $wb.Worksheets.Item("Centre_Registration")
$newSheet.Move($null, $centre)

# Re-fetch a live reference (stale after Move) and clone the standard
# 4-column "module status" header/data block used by every other sheet.
$deb = $wb.Worksheets.Item("Debenture_series_setting")
$centre = $wb.Worksheets.Item("Centre_Registration")
$centre.Range("A1:D2").Copy($deb.Range("A1"))
$deb.Rows.Item(1).RowHeight = 45
$deb.Rows.Item(2).RowHeight = 45
$deb.Range("A2").Value = "Debenture_series_setting"

# ---------------------------------------------------------------------
# 2) InvestmentOpn (now Completed) - replace the stray 19/20/21-digit
#    number columns with the real captured field data.
# ---------------------------------------------------------------------
$inv = $wb.Worksheets.Item("InvestmentOpn")

# Drop the old E:G header/data (19digitNum/20digitNum/21digitNum + huge numbers)
$inv.Range("E1:G2").Clear()

# Clone the header style (bold/fill/wrap) from A1 across the new E1:N1 headers
$inv.Range("A1").Copy() | Out-Null
$inv.Range("E1:N1").PasteSpecial(-4122) | Out-Null

$inv.Range("E1").Value = "accAtBranch"
$inv.Range("F1").Value = "openAmt"
$inv.Range("G1").Value = "IntGLhead"
$inv.Range("H1").Value = "RceiptNo"
$inv.Range("I1").Value = "depstName"
$inv.Range("J1").Value = "reMarks"
$inv.Range("K1").Value = "EnterrelatioN"
$inv.Range("L1").Value = "RemarksTOenter"
$inv.Range("M1").Value = "SPiinfo"
$inv.Range("N1").Value = "RceiptNo"

$inv.Range("A2").Value = "Investment_Opening"
$inv.Range("E2").Value = 0
$inv.Range("F2").Value = 200
$inv.Range("F2").NumberFormat = "General"
$inv.Range("G2").Value = 160
$inv.Range("H2").Value = 3232
$inv.Range("I2").Value = "shradda"
$inv.Range("J2").Value = "No remarks"
$inv.Range("K2").Value = "Daughter"
$inv.Range("L2").Value = "NO"

# Widen the new data columns
$inv.Columns.Item(7).ColumnWidth = 26.8
$inv.Columns.Item(8).ColumnWidth = 14.3

# ---------------------------------------------------------------------
# 3) Update selections. Order matters: the sheet whose Range.Select()
#    runs last becomes the active tab, so InvestmentOpn must be last.
# ---------------------------------------------------------------------
$deb = $wb.Worksheets.Item("Debenture_series_setting")
$deb.Range("K11").Select() | Out-Null

$centre = $wb.Worksheets.Item("Centre_Registration")
$centre.Range("E1:E2").Select() | Out-Null

$inv = $wb.Worksheets.Item("InvestmentOpn")
$inv.Range("O3").Select() | Out-Null

Write-Output "done"
